# Update the cryptos worksheet with the latest scraped price/volume data.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Rows 49/50 also swap
# (Stellar <-> ThetaToken order changed in the source feed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.223.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.348.81'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -5.46%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.342.77'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.26%  '
$ws.Range("E10").Value = '  -8.73%  '
$ws.Range("E11").Value = '  -6.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.52'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -8.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000266'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.64'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.876.90'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '604.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.14'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.264.74'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.336.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.82%  '
$ws.Range("E20").Value = '  -3.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.909'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.93'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '100.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.69%  '
$ws.Range("E26").Value = '  -7.22%  '
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("E28").Value = '  -8.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.72'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.37'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.83%  '
$ws.Range("E32").Value = '  -7.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -14.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.04'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.835.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.34%  '
$ws.Range("E36").Value = '  -5.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '532.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '57.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.61%  '
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.42'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.77%  '
$ws.Range("E41").Value = '  -12.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.67'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.33%  '
$ws.Range("E43").Value = '  -6.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.342'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.78%  '
$ws.Range("E45").Value = '  -6.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.16'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +18.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0414'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.15'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.95%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.42%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.130'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.37%  '
